$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This workbook reorders the 14 data rows (rows 2-15) of the Papaya price
# sheet. Columns A,B,C,E,F,G,H,I,J,K,R are identical across every row, so
# only the varying columns (D,L,M,N,O,P,Q,S,T) are rewritten below, using a
# snapshot of the original row values captured first (so writes never read
# back an already-overwritten cell).

$orig = @{}
$orig[2] = @{ 'D' = 44904; 'L' = 'Primera'; 'M' = 45; 'N' = 15000; 'O' = 15000; 'P' = 15000; 'Q' = '$/bandeja 10 kilos'; 'S' = 1500; 'T' = 10 }
$orig[3] = @{ 'D' = 44904; 'L' = 'Segunda'; 'M' = 60; 'N' = 10000; 'O' = 10000; 'P' = 10000; 'Q' = '$/bandeja 10 kilos'; 'S' = 1000; 'T' = 10 }
$orig[4] = @{ 'D' = 44391; 'L' = 'Primera'; 'M' = 15; 'N' = 1500; 'O' = 1500; 'P' = 1500; 'Q' = '$/kilo (en caja de 15 kilos)'; 'S' = 1500; 'T' = 1 }
$orig[5] = @{ 'D' = 44391; 'L' = 'Segunda'; 'M' = 20; 'N' = 1000; 'O' = 1000; 'P' = 1000; 'Q' = '$/kilo (en caja de 15 kilos)'; 'S' = 1000; 'T' = 1 }
$orig[6] = @{ 'D' = 44195; 'L' = 'Primera'; 'M' = 20; 'N' = 15000; 'O' = 15000; 'P' = 15000; 'Q' = '$/bandeja 10 kilos'; 'S' = 1500; 'T' = 10 }
$orig[7] = @{ 'D' = 44400; 'L' = 'Primera'; 'M' = 25; 'N' = 1500; 'O' = 1500; 'P' = 1500; 'Q' = '$/kilo (en caja de 15 kilos)'; 'S' = 1500; 'T' = 1 }
$orig[8] = @{ 'D' = 44336; 'L' = 'Primera'; 'M' = 10; 'N' = 1500; 'O' = 1500; 'P' = 1500; 'Q' = '$/kilo (en caja de 15 kilos)'; 'S' = 1500; 'T' = 1 }
$orig[9] = @{ 'D' = 44880; 'L' = 'Primera'; 'M' = 200; 'N' = 20000; 'O' = 20000; 'P' = 20000; 'Q' = '$/bandeja 10 kilos'; 'S' = 2000; 'T' = 10 }
$orig[10] = @{ 'D' = 44880; 'L' = 'Segunda'; 'M' = 180; 'N' = 15000; 'O' = 15000; 'P' = 15000; 'Q' = '$/bandeja 10 kilos'; 'S' = 1500; 'T' = 10 }
$orig[11] = @{ 'D' = 44343; 'L' = 'Primera'; 'M' = 20; 'N' = 1700; 'O' = 1700; 'P' = 1700; 'Q' = '$/kilo (en caja de 15 kilos)'; 'S' = 1700; 'T' = 1 }
$orig[12] = @{ 'D' = 44292; 'L' = 'Primera'; 'M' = 50; 'N' = 14000; 'O' = 14000; 'P' = 14000; 'Q' = '$/bandeja 10 kilos'; 'S' = 1400; 'T' = 10 }
$orig[13] = @{ 'D' = 44309; 'L' = 'Primera'; 'M' = 10; 'N' = 1600; 'O' = 1600; 'P' = 1600; 'Q' = '$/kilo (en caja de 15 kilos)'; 'S' = 1600; 'T' = 1 }
$orig[14] = @{ 'D' = 44371; 'L' = 'Primera'; 'M' = 20; 'N' = 1800; 'O' = 1800; 'P' = 1800; 'Q' = '$/kilo (en caja de 15 kilos)'; 'S' = 1800; 'T' = 1 }
$orig[15] = @{ 'D' = 44371; 'L' = 'Segunda'; 'M' = 30; 'N' = 1200; 'O' = 1200; 'P' = 1200; 'Q' = '$/kilo (en caja de 15 kilos)'; 'S' = 1200; 'T' = 1 }

# target row -> source row (original row whose values now occupy it)
$rowMap = @{}
$rowMap[2] = 4
$rowMap[3] = 5
$rowMap[4] = 2
$rowMap[5] = 3
$rowMap[6] = 14
$rowMap[7] = 15
$rowMap[8] = 6
$rowMap[9] = 11
$rowMap[10] = 13
$rowMap[11] = 9
$rowMap[12] = 10
$rowMap[13] = 8
$rowMap[14] = 12
$rowMap[15] = 7

foreach ($targetRow in 2..15) {
    $srcRow = $rowMap[$targetRow]
    $vals = $orig[$srcRow]
    $ws.Cells.Item($targetRow, 4).Value = $vals['D']
    $ws.Cells.Item($targetRow, 12).Value = $vals['L']
    $ws.Cells.Item($targetRow, 13).Value = $vals['M']
    $ws.Cells.Item($targetRow, 14).Value = $vals['N']
    $ws.Cells.Item($targetRow, 15).Value = $vals['O']
    $ws.Cells.Item($targetRow, 16).Value = $vals['P']
    $ws.Cells.Item($targetRow, 17).Value = $vals['Q']
    $ws.Cells.Item($targetRow, 19).Value = $vals['S']
    $ws.Cells.Item($targetRow, 20).Value = $vals['T']
}

"done"
